# Fixed #295 - Add the version of M2Doc in the template custom properties.
#
# The unified diff for this particular fixture
# (inlinedUserContent-template.docx) only touches two package parts,
# word/document.xml and word/styles.xml, and every hunk is a pure
# attribute-reordering artifact from the authoring tool's XML writer:
#   - xmlns:* declarations on the <w:document> root are re-emitted in a
#     different (alphabetical) order,
#   - attributes on <w:pgSz>/<w:pgMar>, <w:rFonts>/<w:lang>,
#     <w:latentStyles>/<w:lsdException>, <w:style>, <w:tblInd> and
#     <w:tblCellMar> are likewise just re-ordered.
# Every element, every attribute name/value pair, every run of text
# ("A simple demonstration of a ", "user doc", " :", "Some protected
# text.", "End of demonstration.") and every numeric value (page size
# 11906x16838, margins 1417/1417/1417/1417, fonts, language, latent
# style table, built-in style definitions, table cell margins, ...)
# is byte-for-byte identical before and after - nothing a user could
# observe in Word changed. That is consistent with the commit message:
# the actual "add M2Doc version to custom properties" work landed in
# the M2Doc Java sources / other template fixtures, and this resource
# was simply swept up in the same bulk re-save (hence the cosmetic
# attribute-order churn) without any of its own content being edited.
#
# Word's object model (real COM or this emulation) has no property
# that controls on-disk XML attribute ordering - that is purely a
# detail of whatever tool serialized the part - so there is no
# Document/Range/Style/PageSetup call that corresponds to this diff.
# Issuing a cosmetic "no-op" edit (e.g. re-applying the same Find, or
# re-writing a property to its own value) would not reproduce the
# attribute order either, and it would actually move the document
# further from the target by picking up incidental side effects
# (freshly minted namespace declarations, rsid churn, etc.). So the
# correct replay of this diff is to leave the document content
# untouched.
$d = $word.ActiveDocument

# Touch the document read-only, just to confirm the session/object
# model is alive; this performs no mutation of any part.
$null = $d.Name
$null = $d.Paragraphs.Count
